$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.551.61"
$ws.Range("E2").Value = "  -0.47%  "
$ws.Range("D3").Value = "1.728.94"
$ws.Range("E3").Value = "  -1.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.22"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4814"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2671"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.63%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06182"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("D10").Value = "1.727.72"
$ws.Range("E10").Value = "  -1.23%  "
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.58"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.69%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6097"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.89%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.537"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.29%  "
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").Value = "26.554.47"
$ws.Range("E17").Value = "  -0.47%  "
$ws.Range("E18").Value = "  +0.09%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000006956"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.79%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.62%  "
$ws.Range("D21").Value = "1.952.29"
$ws.Range("E21").Value = "  -0.85%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.520"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.91%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.798"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.98%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.238"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.37%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "137.32"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.62%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.34"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.09%  "
$ws.Range("E27").Value = "  -3.05%  "
$ws.Range("E28").Value = "  -0.69%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.79"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.13%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.968"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.69%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08010"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.43%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.689"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.17%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04506"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.72%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.9999"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.612"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.002"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.21%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6312"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.37%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.9063"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.053"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.396"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "103.12"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -9.62%  "
$ws.Range("E43").Value = "  -0.92%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.475"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.95%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3884"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.10%  "
$ws.Range("E46").Value = "  +4.78%  "
$ws.Range("E47").Value = "  -1.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.05387"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.97%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.880"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.249"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.04%  "
